# Updated symbol list on Sat Jan 21 21:41:46 UTC 2023 with GitHub Actions
# Refresh the Price (D) / Volume(1h) (E) columns with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.30%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'35.58"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'9.99%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.070"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.30%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07817"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.18%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'-1.27%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.095"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.05%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.048"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'5.93%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9296"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.75%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.09422"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-5.37%"
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'3.53%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08558"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.49%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03493"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'5.71%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09959"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.97%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001490"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.67%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005727"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.01%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.478"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.73%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.068"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-5.35%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3406"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.50%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1322"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-1.11%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.540"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.72%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'7.17%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04669"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.96%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001232"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.39%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004540"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'3.93%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.66%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-19.94%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01775"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'4.03%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04713"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.65%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007986"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.39%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1421"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.03%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008003"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-18.09%"
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'6.91%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009082"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-6.39%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006206"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.27%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.76%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'4.064"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'45.44%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002693"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'35.52%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.76%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.76%"
$ws.Range("E51").Style = "Normal"
